# Update weather_data.xlsx sensor readings (rows 2-11) with refreshed values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rainNone = "Дождь - 0.0`nЛивень - 0.0`nСнег - 0.0"
$rainLight = "Дождь - 0.8`nЛивень - 0.0`nСнег - 0.07"

$data = @(
    @{ Row = 2;  F = 45580.99404122325; E = $rainNone },
    @{ Row = 3;  F = 45580.99402684862; E = $rainNone },
    @{ Row = 4;  F = 45580.99402511361; E = $rainLight },
    @{ Row = 5;  F = 45580.99401505614; E = $rainNone },
    @{ Row = 6;  F = 45580.99400483284; E = $rainNone },
    @{ Row = 7;  F = 45580.99399292153; E = $rainNone },
    @{ Row = 8;  F = 45580.99398033474; E = $rainNone },
    @{ Row = 9;  F = 45580.99396962391; E = $rainNone },
    @{ Row = 10; F = 45580.99390919802; E = $rainLight },
    @{ Row = 11; F = 45580.99379513197; E = $rainLight }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = 5.9
    $ws.Cells.Item($r, 2).Value = 6
    $ws.Cells.Item($r, 3).Value = "ЮЗ"
    $ws.Cells.Item($r, 4).Value = 994.9
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
    # Writing a multi-line value stamps an explicit row height; restore
    # the row to its natural (non-custom) height to match the source edit.
    $ws.Rows.Item($r).EntireRow.AutoFit()
}
